$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "002"

$ws.Range("N2").Value = "2020-06-30 00:00:00"

$ws.Range("O2").Value = 52524058.99
$ws.Range("P2").Value = 248.9172347745
$ws.Range("Q2").Value = 266182971.82
$ws.Range("R2").Value = 1261.4700874909
$ws.Range("S2").Value = 22534272.56
$ws.Range("T2").Value = 106.7923713656
$ws.Range("U2").Value = -27112290.21
$ws.Range("V2").Value = -128.4880954985
$ws.Range("Y2").Value = 27112290.21
$ws.Range("Z2").Value = 128.4880954985
$ws.Range("AA2").Value = -6035650.01
$ws.Range("AB2").Value = -28.6036026051
$ws.Range("AC2").Value = 21101013.37
$ws.Range("AD2").Value = 142.7405266799
